# Update "bibi" retention metrics (2024 period_index=1 row and 2025 period_index=0 row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: cohort_year 2024, period_index 1 -> num_customers 109 -> 111
$ws.Range("C36").Value = 111
# retention_rate = num_customers / cohort_size
$ws.Range("E36").Value = 111 / 1930

# Row 37: cohort_year 2025, period_index 0 -> num_customers / cohort_size 664 -> 667
$ws.Range("C37").Value = 667
$ws.Range("D37").Value = 667
$ws.Range("E37").Value = 1
